$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.279.93'
$ws.Range('E2').Value = '  +0.05%  '
$ws.Range('D3').Value = '3.494.92'
$ws.Range('E3').Value = '  +0.09%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '''588.93'
$ws.Range('E5').Value = '  +0.34%  '
$ws.Range('D6').Value = '''133.61'
$ws.Range('E6').Value = '  -0.34%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '''0.485'
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E9').Value = '  +6.35%  '
$ws.Range('E10').Value = '  +0.06%  '
$ws.Range('D11').Value = '''0.387'
$ws.Range('E11').Value = '  +3.03%  '
$ws.Range('D12').Value = '4.092.16'
$ws.Range('E12').Value = '  +0.18%  '
$ws.Range('E13').Value = '  +0.71%  '
$ws.Range('E14').Value = '  -0.22%  '
$ws.Range('D15').Value = '3.501.01'
$ws.Range('E15').Value = '  +0.26%  '
$ws.Range('D16').Value = '64.249.06'
$ws.Range('E16').Value = '  -0.07%  '
$ws.Range('D17').Value = '''25.28'
$ws.Range('E17').Value = '  +0.11%  '
$ws.Range('D18').Value = '''10.04'
$ws.Range('D19').Value = '''5.77'
$ws.Range('E19').Value = '  +0.56%  '
$ws.Range('D20').Value = '''13.52'
$ws.Range('E20').Value = '  -0.56%  '
$ws.Range('D21').Value = '''386.70'
$ws.Range('E21').Value = '  -0.25%  '
$ws.Range('D22').Value = '''0.580'
$ws.Range('E22').Value = '  +2.67%  '
$ws.Range('D23').Value = '3.635.25'
$ws.Range('E23').Value = '  +0.13%  '
$ws.Range('D24').Value = '''74.11'
$ws.Range('E24').Value = '  -0.30%  '
$ws.Range('E25').Value = '  +0.12%  '
$ws.Range('E26').Value = '  -0.47%  '
$ws.Range('E27').Value = '  +2.27%  '
$ws.Range('D28').Value = '''0.999'
$ws.Range('E28').Value = '  -0.08%  '
$ws.Range('D29').Value = '''7.36'
$ws.Range('E29').Value = '  -0.30%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').Value = '''2.25'
$ws.Range('E30').Value = '  +0.80%  '
$ws.Range('B31').Value = 'Fetch.AI'
$ws.Range('C31').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D31').Value = '''1.50'
$ws.Range('E31').Value = '  -1.39%  '
$ws.Range('D32').Value = '''8.15'
$ws.Range('E32').Value = '  -0.88%  '
$ws.Range('E33').Value = '  +3.53%  '
$ws.Range('D34').Value = '3.524.83'
$ws.Range('E34').Value = '  +0.36%  '
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('D36').Value = '''23.29'
$ws.Range('E36').Value = '  -0.68%  '
$ws.Range('E37').Value = '  +1.73%  '
$ws.Range('D38').Value = '''6.91'
$ws.Range('E39').Value = '  +0.21%  '
$ws.Range('D40').Value = '''165.47'
$ws.Range('E40').Value = '  +1.88%  '
$ws.Range('D41').Value = '''0.0785'
$ws.Range('E41').Value = '  +0.51%  '
$ws.Range('E42').Value = '  +0.32%  '
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('D45').Value = '''24.41'
$ws.Range('E45').Value = '  -4.20%  '
$ws.Range('D46').Value = '''1.17'
$ws.Range('E46').Value = '  -0.29%  '
$ws.Range('E47').Value = '  -1.12%  '
$ws.Range('D48').Value = '2.423.91'
$ws.Range('E48').Value = '  -1.87%  '
$ws.Range('D49').Value = '''6.82'
$ws.Range('E49').Value = '  +1.04%  '
$ws.Range('D50').Value = '''0.917'
$ws.Range('E50').Value = '  +1.70%  '
$ws.Range('E51').Value = '  -0.60%  '
